$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Future Improvements" bullet list (numId=6): the second bullet item was
#    left empty in the original document; fill it in with the new text about
#    data-structure improvements.
# ---------------------------------------------------------------------------
$added = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (-not $added -and $t.Trim() -eq "" -and $p.Range.ListFormat.ListType -eq 3) {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text -like "*Various heuristics can be used*") {
            $p.Range.Text = "Data structures and management " + [char]0x2013 + " there is a lot of potential to remove redundancy from data structures such as indexing instead of using strings, using trees to store data that need to be constantly sorted and eliminating redundant copying."
            $added = $true
        }
    }
}

# ---------------------------------------------------------------------------
# 2. "Open Issues" bullet list (numId=3): the last two bullets were
#       "Fix dependency generation in tests"
#       <empty>
#    and become a single bullet reading "Dependency violation issues".
#    First fold the trailing empty bullet into the "Fix dependency..." one by
#    deleting the paragraph mark that separates them, then update the text.
# ---------------------------------------------------------------------------
$fix = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Fix dependency generation in tests") {
        $fix = $p
    }
}
if ($fix -ne $null) {
    $nxt = $fix.Next()
    if ($nxt -ne $null -and $nxt.Range.Text.Trim() -eq "") {
        $markStart = $fix.Range.End - 1
        $markEnd = $fix.Range.End
        $d.Range($markStart, $markEnd).Delete()
    }
}

$d.Content.Find.Execute("Fix dependency generation in tests", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dependency violation issues", 2)
